# Swap the contents of columns C ("codeforiati:group-name") and
# D ("codeforiati:group-code") so that column C holds the group code
# and column D holds the group name, for the header row and every
# data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
